$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.213.99"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").Value = "1.859.14"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "236.76"
$ws.Range("E5").Value = "  +1.06%  "

# Row 6
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("D7").Value = "0.4672"
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").Value = "0.2863"
$ws.Range("E8").Value = "  +0.96%  "

# Row 9
$ws.Range("D9").Value = "0.06527"
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").Value = "21.80"
$ws.Range("E10").Value = "  +4.72%  "

# Row 11
$ws.Range("D11").Value = "0.07936"
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").Value = "97.04"
$ws.Range("E12").Value = "  +0.81%  "

# Row 13
$ws.Range("D13").Value = "1.867.72"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14
$ws.Range("D14").Value = "5.165"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15
$ws.Range("D15").Value = "0.6794"
$ws.Range("E15").Value = "  +1.36%  "

# Row 16
$ws.Range("D16").Value = "266.38"
$ws.Range("E16").Value = "  -4.51%  "

# Row 17
$ws.Range("D17").Value = "30.197.98"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").Value = "13.64"
$ws.Range("E18").Value = "  +7.91%  "

# Row 19
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.10%  "

# Row 20
$ws.Range("D20").Value = "0.000007365"
$ws.Range("E20").Value = "  +1.62%  "

# Row 21
$ws.Range("D21").Value = "2.113.73"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  -2.93%  "

# Row 23
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").Value = "6.190"
$ws.Range("E24").Value = "  +0.65%  "

# Row 25
$ws.Range("D25").Value = "167.00"
$ws.Range("E25").Value = "  +1.48%  "

# Row 26
$ws.Range("D26").Value = "9.204"
$ws.Range("E26").Value = "  -1.01%  "

# Row 27
$ws.Range("D27").Value = "18.82"
$ws.Range("E27").Value = "  -0.63%  "

# Row 28
$ws.Range("D28").Value = "1.952"
$ws.Range("E28").Value = "  +2.58%  "

# Row 29
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").Value = "  +2.55%  "

# Row 30
$ws.Range("D30").Value = "0.09861"
$ws.Range("E30").Value = "  +3.01%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.366"
$ws.Range("E31").Value = "  -0.76%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.479"
$ws.Range("E32").Value = "  +0.77%  "

# Row 33
$ws.Range("D33").Value = "4.051"
$ws.Range("E33").Value = "  -1.24%  "

# Row 34
$ws.Range("D34").Value = "0.04705"
$ws.Range("E34").Value = "  +0.22%  "

# Row 35
$ws.Range("D35").Value = "1.127"
$ws.Range("E35").Value = "  +2.70%  "

# Row 36
$ws.Range("D36").Value = "0.7011"
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
$ws.Range("D38").Value = "0.01872"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("E39").Value = "  +4.12%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.229"
$ws.Range("E40").Value = "  -1.58%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "74.48"
$ws.Range("E41").Value = "  +2.09%  "

# Row 42
$ws.Range("D42").Value = "1.936"
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
$ws.Range("D43").Value = "0.8486"
$ws.Range("E43").Value = "  +0.33%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.13%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4156"
$ws.Range("E45").Value = "  -0.26%  "

# Row 46
$ws.Range("D46").Value = "103.36"
$ws.Range("E46").Value = "  -0.51%  "

# Row 47
$ws.Range("D47").Value = "955.62"
$ws.Range("E47").Value = "  -3.06%  "

# Row 48
$ws.Range("D48").Value = "7.144"
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("D49").Value = "9.181"
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "34.09"
$ws.Range("E50").Value = "  +0.37%  "

# Row 51
$ws.Range("D51").Value = "0.05647"
$ws.Range("E51").Value = "  +0.59%  "
